$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Ready for handoff" -> "In Translation" (Status column on every sheet)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"

# ---------------------------------------------------------------------------
# 2) Narrow the "Latest HO Xliff Generate Date"-style status columns
#    (Overview!E:F and the "Status" column on the language sheets)
# ---------------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 12.5   # column E
$wsOverview.Columns.Item(6).ColumnWidth = 12.5   # column F

$wsZhCn.Columns.Item(3).ColumnWidth = 12.5        # column C
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5        # column C
